$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 101, shifting the existing rows 101-176 down to 102-177.
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with the new record.
$ws.Range("A101").Value = 4
$ws.Range("B101").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C101").Value = "Los Lagos"
$ws.Range("D101").Value = 45233
$ws.Range("E101").Value = 10
$ws.Range("F101").Value = 100112026
$ws.Range("G101").Value = "Haba"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 120
$ws.Range("K101").Value = 16000
$ws.Range("L101").Value = 16000
$ws.Range("M101").Value = 16000
$ws.Range("N101").Value = "$/saco 25 kilos"
$ws.Range("O101").Value = "Región Metropolitana"
$ws.Range("P101").Value = 640
$ws.Range("Q101").Value = 25
$ws.Range("R101").Value = "Hortaliza"
